# Cover letter rewrite: personal header block added, body replaced with the
# Goldman Sachs Software Engineer version, closing block reformatted to
# match (Cambria font throughout, tightened paragraph spacing on the
# header/signature lines).
$d = $word.ActiveDocument

$newBodyXml = @'
<w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>Deep Manek</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>201-241-5648</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>dpmanek@gmail.com</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>October 7, 2023</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>Dear Hiring Manager,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>I am writing to express my interest in the Software Engineer position at Goldman Sachs, as described in your recent job listing. With a strong academic background in Computer Science and substantial professional experience as a software developer, I believe I am well-suited to contribute effectively to the innovative and client-centric ethos of your organization.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t xml:space="preserve">My </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>Master’s degree in Computer Science</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t xml:space="preserve"> from Stevens Institute of Technology, combined with my hands-on experience at </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>At</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t xml:space="preserve"> Last Sportswear and Larsen &amp; Toubro Infotech, has equipped me with a robust skill set in web programming, software design, and system architecture. In my recent role, I spearheaded the transformation of a major e-commerce website, harnessing ReactJS to effectuate a 30% sales increase. Furthermore, I collaborated across teams to deploy RESTful APIs and optimize data storage solutions, ensuring efficient and scalable applications.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>In addition to my technical prowess, my experience working in diverse environments has honed my communication skills and ability to collaborate with both technical and non-technical teams. I am adept at managing multiple projects, with a keen attention to detail and a proven ability to meet tight deadlines, qualities that I believe align with the role’s requirements. My experience with UI/UX development and API design, coupled with a foundational understanding of the financial industry, makes me a unique candidate who can bridge the gap between technology and its practical application in a financial context.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>Goldman Sachs's commitment to innovation, as well as fostering diversity and inclusion, resonates strongly with me. I am excited about the opportunity to be a part of embedded engineering teams, leveraging the latest technologies to deliver groundbreaking solutions that maintain Goldman Sachs's leading position in the industry.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>I appreciate your consideration and look forward to the opportunity to further discuss how I can be an asset to Goldman Sachs. Please find my attached resume for a detailed review of my credentials. I am available at your earliest convenience for an interview and can be reached at 201-241-5648 or via email at dpmanek@gmail.com.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>Thank you for your time and consideration.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>Sincerely,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
        <w:t>Deep Manek</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    
</w:body>
'@

# Replacing the whole story range (but not the sectPr, which Word keeps)
# in one shot keeps every paragraph's run-level formatting (Cambria rFonts)
# and the proofErr spans exactly as authored.
[void]$d.Content.InsertXML($newBodyXml)

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
